# Refresh the cryptos snapshot: Price (column D) and Volume(1h) (column E),
# matching the scheduled GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store an otherwise-numeric-looking
# string (e.g. '507.28') as literal text, matching every Price cell in
# this sheet -- exactly as if a person typed '507.28 into the cell.
$textPrefix = [string][char]39

$ws.Range('D2').Value = '58.206.74'
$ws.Range('E2').Value = '  -4.64%  '
$ws.Range('D3').Value = '2.574.55'
$ws.Range('E3').Value = '  -3.54%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = $textPrefix + '507.28'
$ws.Range('E5').Value = '  -4.37%  '
$ws.Range('D6').Value = $textPrefix + '144.36'
$ws.Range('E6').Value = '  -7.31%  '
$ws.Range('D8').Value = $textPrefix + '0.571'
$ws.Range('E8').Value = '  -2.12%  '
$ws.Range('D9').Value = '2.584.83'
$ws.Range('E9').Value = '  -3.80%  '
$ws.Range('D10').Value = $textPrefix + '6.28'
$ws.Range('E10').Value = '  -4.03%  '
$ws.Range('E11').Value = '  -4.78%  '
$ws.Range('E12').Value = '  -5.11%  '
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('D14').Value = '3.025.40'
$ws.Range('E14').Value = '  -3.70%  '
$ws.Range('D15').Value = '58.189.71'
$ws.Range('E15').Value = '  -4.72%  '
$ws.Range('D16').Value = $textPrefix + '21.03'
$ws.Range('E16').Value = '  -4.78%  '
$ws.Range('E17').Value = '  -4.20%  '
$ws.Range('D18').Value = '2.588.98'
$ws.Range('E18').Value = '  -2.97%  '
$ws.Range('E19').Value = '  -5.13%  '
$ws.Range('D20').Value = $textPrefix + '341.98'
$ws.Range('E20').Value = '  -3.68%  '
$ws.Range('D21').Value = $textPrefix + '10.30'
$ws.Range('E21').Value = '  -3.87%  '
$ws.Range('E22').Value = '  -4.26%  '
$ws.Range('D23').Value = $textPrefix + '1.00'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').Value = $textPrefix + '60.60'
$ws.Range('E24').Value = '  -2.00%  '
$ws.Range('D25').Value = $textPrefix + '0.418'
$ws.Range('E25').Value = '  -3.30%  '
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('D27').Value = '2.685.67'
$ws.Range('E27').Value = '  -3.86%  '
$ws.Range('E28').Value = '  -5.80%  '
$ws.Range('E29').Value = '  -5.20%  '
$ws.Range('E30').Value = '  -5.19%  '
$ws.Range('D31').Value = $textPrefix + '0.999'
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('D32').Value = $textPrefix + '6.11'
$ws.Range('E32').Value = '  -1.23%  '
$ws.Range('D33').Value = $textPrefix + '18.81'
$ws.Range('E33').Value = '  -3.77%  '
$ws.Range('D34').Value = $textPrefix + '149.06'
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('E35').Value = '  -5.60%  '
$ws.Range('D36').Value = $textPrefix + '0.951'
$ws.Range('E36').Value = '  +6.89%  '
$ws.Range('D37').Value = $textPrefix + '3.98'
$ws.Range('E37').Value = '  -3.74%  '
$ws.Range('E38').Value = '  -5.96%  '
$ws.Range('D39').Value = $textPrefix + '0.849'
$ws.Range('E39').Value = '  -7.89%  '
$ws.Range('D40').Value = $textPrefix + '36.00'
$ws.Range('E40').Value = '  -2.37%  '
$ws.Range('D41').Value = $textPrefix + '288.77'
$ws.Range('E41').Value = '  -5.87%  '
$ws.Range('E42').Value = '  -5.44%  '
$ws.Range('E43').Value = '  -6.86%  '
$ws.Range('D44').Value = $textPrefix + '0.0990'
$ws.Range('E44').Value = '  -2.70%  '
$ws.Range('D45').Value = $textPrefix + '0.995'
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('E46').Value = '  -6.48%  '
$ws.Range('D47').Value = $textPrefix + '19.28'
$ws.Range('E47').Value = '  -6.44%  '
$ws.Range('D48').Value = $textPrefix + '0.0535'
$ws.Range('E48').Value = '  -5.48%  '
$ws.Range('E49').Value = '  -1.16%  '
$ws.Range('E50').Value = '  -5.71%  '
$ws.Range('E51').Value = '  -7.21%  '

Write-Output "applied 79 cell updates"
